{"js": "// \"\u00c6ndret vare til produkt\" \u2014 replace the Danish word \"vare\" (\"item\") with\n// \"produkt\" (\"product\") in the two user-test observation paragraphs that use\n// it as a common noun. The occurrence inside \"svaret\" (=\"the answer\") must be\n// left untouched, so we match whole phrases rather than doing a blind\n// substring replace.\nconst body = context.document.body;\n\n// 1) \"...s\u00f8g efter en vare og taster...\" -> \"...s\u00f8g efter et produkt og taster...\"\n//    (the indefinite article also changes because \"produkt\" is neuter in\n//    Danish while \"vare\" is common gender).\nconst hit1 = body.search(\"en vare og taster\", { matchCase: true, matchWholeWord: false });\nhit1.load(\"text\");\n\n// 2) \"...s\u00f8g efter vare og klikker...\" -> \"...s\u00f8g efter produkt og klikker...\"\nconst hit2 = body.search(\"efter vare og klikker\", { matchCase: true, matchWholeWord: false });\nhit2.load(\"text\");\n\nawait context.sync();\n\nif (hit1.items.length !== 1) {\n  throw new Error(`expected exactly 1 match for \"en vare og taster\", found ${hit1.items.length}`);\n}\nif (hit2.items.length !== 1) {\n  throw new Error(`expected exactly 1 match for \"efter vare og klikker\", found ${hit2.items.length}`);\n}\n\nhit1.items[0].insertText(\"et produkt og taster\", \"Replace\");\nhit2.items[0].insertText(\"efter produkt og klikker\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# \"\u00c6ndret vare til produkt\" - replace the Danish word \"vare\" (\"item\") with\n# \"produkt\" (\"product\") in the two user-test observation paragraphs that use\n# it as a common noun. The occurrence inside \"svaret\" (=\"the answer\") must be\n# left untouched, so we search/replace whole phrases rather than doing a\n# blind substring replace of \"vare\".\n\n$d = $word.ActiveDocument\n\n# 1) \"...s\u00f8g efter en vare og taster...\" -> \"...s\u00f8g efter et produkt og taster...\"\n#    (the indefinite article also changes because \"produkt\" is neuter in\n#    Danish while \"vare\" is common gender).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$found1 = $find1.Execute(\n    \"en vare og taster\",   # FindText\n    $false,                # MatchCase\n    $false,                # MatchWholeWord\n    $false,                # MatchWildcards\n    $false,                # MatchSoundsLike\n    $false,                # MatchAllWordForms\n    $true,                 # Forward\n    1,                     # Wrap (wdFindContinue)\n    $false,                # Format\n    \"et produkt og taster\",# ReplaceWith\n    2                      # Replace (wdReplaceAll)\n)\nif (-not $found1) {\n    throw \"Could not find 'en vare og taster' to replace.\"\n}\n\n# 2) \"...s\u00f8g efter vare og klikker...\" -> \"...s\u00f8g efter produkt og klikker...\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute(\n    \"efter vare og klikker\",    # FindText\n    $false,                     # MatchCase\n    $false,                     # MatchWholeWord\n    $false,                     # MatchWildcards\n    $false,                     # MatchSoundsLike\n    $false,                     # MatchAllWordForms\n    $true,                      # Forward\n    1,                          # Wrap (wdFindContinue)\n    $false,                     # Format\n    \"efter produkt og klikker\", # ReplaceWith\n    2                           # Replace (wdReplaceAll)\n)\nif (-not $found2) {\n    throw \"Could not find 'efter vare og klikker' to replace.\"\n}\n"}
